$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells that would otherwise be parsed as numbers
$priceTextCells = @("D5", "D6", "D8", "D9", "D10", "D13", "D16", "D18", "D19", "D21", "D23", "D24", "D26", "D27", "D30", "D31", "D32", "D34", "D36", "D37", "D38", "D39", "D40", "D41", "D44", "D45", "D47", "D49")
foreach ($addr in $priceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "34.555.23"
$ws.Range("E2").Value = "  -0.07%  "

$ws.Range("D3").Value = "1.812.59"
$ws.Range("E3").Value = "  +0.99%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "228.67"
$ws.Range("E5").Value = "  +0.58%  "

$ws.Range("D6").Value = "0.578"
$ws.Range("E6").Value = "  +3.98%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").Value = "34.96"
$ws.Range("E8").Value = "  +6.51%  "

$ws.Range("D9").Value = "0.304"
$ws.Range("E9").Value = "  +2.39%  "

$ws.Range("D10").Value = "0.0698"
$ws.Range("E10").Value = "  +0.74%  "

$ws.Range("E11").Value = "  +1.02%  "

$ws.Range("D12").Value = "2.075.30"
$ws.Range("E12").Value = "  +0.99%  "

$ws.Range("D13").Value = "11.32"
$ws.Range("E13").Value = "  +1.03%  "

$ws.Range("D14").Value = "1.822.85"
$ws.Range("E14").Value = "  +1.49%  "

$ws.Range("E15").Value = "  +1.83%  "

$ws.Range("D16").Value = "4.56"
$ws.Range("E16").Value = "  +5.12%  "

$ws.Range("D17").Value = "34.540.07"
$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("D18").Value = "69.34"
$ws.Range("E18").Value = "  +0.57%  "

$ws.Range("D19").Value = "247.29"
$ws.Range("E19").Value = "  +0.40%  "

$ws.Range("D20").Value = "0.0₃0801"
$ws.Range("E20").Value = "  -0.12%  "

$ws.Range("D21").Value = "11.51"
$ws.Range("E21").Value = "  +0.95%  "

$ws.Range("E22").Value = "  +0.06%  "

$ws.Range("D23").Value = "4.21"
$ws.Range("E23").Value = "  +0.84%  "

$ws.Range("D24").Value = "173.39"
$ws.Range("E24").Value = "  -0.95%  "

$ws.Range("E25").Value = "  +2.42%  "

$ws.Range("D26").Value = "8.03"
$ws.Range("E26").Value = "  +9.09%  "

$ws.Range("D27").Value = "16.85"
$ws.Range("E27").Value = "  +0.91%  "

$ws.Range("E28").Value = "  +3.15%  "

$ws.Range("E29").Value = "  -0.09%  "

$ws.Range("D30").Value = "4.08"
$ws.Range("E30").Value = "  +1.25%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "3.93"
$ws.Range("E31").Value = "  +3.27%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "0.0538"
$ws.Range("E32").Value = "  +2.45%  "

$ws.Range("E33").Value = "  +0.83%  "

$ws.Range("D34").Value = "1.85"
$ws.Range("E34").Value = "  +1.05%  "

$ws.Range("D35").Value = "1.398.48"
$ws.Range("E35").Value = "  -1.82%  "

$ws.Range("D36").Value = "0.680"
$ws.Range("E36").Value = "  -0.05%  "

$ws.Range("D37").Value = "2.49"
$ws.Range("E37").Value = "  -2.40%  "

$ws.Range("D38").Value = "1.07"
$ws.Range("E38").Value = "  -0.08%  "

$ws.Range("D39").Value = "0.0192"
$ws.Range("E39").Value = "  +0.50%  "

$ws.Range("D40").Value = "84.08"
$ws.Range("E40").Value = "  -0.64%  "

$ws.Range("D41").Value = "0.976"
$ws.Range("E41").Value = "  +2.94%  "

$ws.Range("E42").Value = "  +3.17%  "

$ws.Range("E43").Value = "  +0.10%  "

$ws.Range("D44").Value = "1.15"
$ws.Range("E44").Value = "  +6.06%  "

$ws.Range("D45").Value = "13.19"
$ws.Range("E45").Value = "  -5.51%  "

$ws.Range("E46").Value = "  -2.49%  "

$ws.Range("D47").Value = "6.03"
$ws.Range("E47").Value = "  -1.34%  "

$ws.Range("D48").Value = "1.973.92"
$ws.Range("E48").Value = "  +0.91%  "

$ws.Range("D49").Value = "105.48"
$ws.Range("E49").Value = "  -0.11%  "

$ws.Range("E50").Value = "  +1.02%  "

$ws.Range("E51").Value = "  +0.05%  "
